$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the meme text for row 26 (sticker id in B26 stays the same)
$ws.Range("A26").Value = "esos no son"

# 2. Append the newly-linked sticker groups as three new Meme/StickerID rows
$ws.Range("A144").Value = "y ahora | y ahora? | ¿y ahora?"
$ws.Range("B144").Value = "CAACAgEAAxkBAAIiQGA3MJ0T6oZfCnmpDt0Wo_BLNPzQAAIgAAPFhdoNozL7WHG6afYeBA | CAACAgEAAxkBAAIiQmA3MKbhpXKSWxPafBjIFNl6BQTEAAIhAAPFhdoNhf0gTTgqUUAeBA | CAACAgEAAxkBAAIiRGA3MKwks23eTHdHugbWPestq1zZAAIiAAPFhdoNJ7DkTKsET7geBA"

$ws.Range("A145").Value = "me vale verga*"
$ws.Range("B145").Value = "CAACAgEAAxkBAAIiTWA3MkPxDdyR6te5uTzsVg7aDpxPAALHAAPFhdoNMefIKi5yx-QeBA | CAACAgEAAxkBAAIiT2A3MkXoi3GN7KYKb0J5ZXDKigqpAALIAAPFhdoNp7j9ONRjbrgeBA | CAACAgEAAxkBAAIiUWA3MkhF1mHkiQWQaNe2syhd0TcXAALJAAPFhdoNzHqODZBNH1AeBA"

$ws.Range("A146").Value = "i need you"
$ws.Range("B146").Value = "CAACAgEAAxkBAAIiaGA3OfbjMHG_-IBUjOk7hS5Hi7ofAALUAAPFhdoNLdWWYKY6nYYeBA | CAACAgEAAxkBAAIiamA3OfhWMf_uMMADCu-AueU0i4J2AALTAAPFhdoNpbNkS3ptGZYeBA | CAACAgEAAxkBAAIibGA3Ofoh026trXupxeNEB-JXfdFWAALSAAPFhdoN9FkLZizXt8weBA"

# The three new rows carry their own (slightly larger / unbolded) Arial 11 look
$newRows = $ws.Range("A144:B146")
$newRows.Font.Name = "Arial"
$newRows.Font.Size = 11
$newRows.Font.ThemeColor = 1

# 3. Drop the stray trailing blank row at the end of the sheet
$ws.Rows.Item(999).Delete()
